# fall 22 week 16 day after inputs
# Adds a new "Week 50" column (AY) to the InningCounts sheet and fills in
# the inning counts that were recorded for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the week 50 column
$ws.Range("AY1").Value = "Week 50"

# New inning-count values recorded for week 50
$ws.Range("AY4").Value = 7
$ws.Range("AY6").Value = 5.5
$ws.Range("AY8").Value = 10
$ws.Range("AY9").Value = 3
$ws.Range("AY10").Value = 4.5

# Leave the selection where Excel would land after typing the last entry
$ws.Range("AY14").Select() | Out-Null
